$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: new row R should receive the D/J/K/L/M/P block that currently
# (before this edit) lives in row $mapping[R]. This is the permutation of
# weekly price records described by the commit ("Fruta / hortaliza, semanal").
$mapping = @{}
$mapping[2]  = 40
$mapping[3]  = 9
$mapping[4]  = 5
$mapping[5]  = 34
$mapping[6]  = 10
$mapping[7]  = 13
$mapping[8]  = 22
$mapping[9]  = 4
$mapping[10] = 30
$mapping[11] = 14
$mapping[12] = 42
$mapping[13] = 25
$mapping[14] = 41
$mapping[15] = 21
$mapping[16] = 19
$mapping[17] = 43
$mapping[18] = 39
$mapping[19] = 8
$mapping[20] = 3
$mapping[21] = 31
$mapping[22] = 33
$mapping[23] = 38
$mapping[24] = 26
$mapping[25] = 20
$mapping[26] = 28
$mapping[27] = 16
$mapping[28] = 12
$mapping[29] = 37
$mapping[30] = 36
$mapping[31] = 24
$mapping[32] = 18
$mapping[33] = 44
$mapping[34] = 32
$mapping[35] = 11
$mapping[36] = 27
$mapping[37] = 46
$mapping[38] = 23
$mapping[39] = 7
$mapping[40] = 2
$mapping[41] = 29
$mapping[42] = 45
$mapping[43] = 6
$mapping[44] = 17
$mapping[45] = 15
$mapping[46] = 35

# Snapshot the current D,J,K,L,M,P values for every data row (2-46) before
# any writes happen, so the permutation can be applied safely even though
# it is not a simple pairwise swap. Value2 is used (rather than Value) since
# it reliably round-trips the raw numeric/date serial for these cells.
$orig = @{}
foreach ($r in 2..46) {
    $orig[$r] = @(
        $ws.Cells.Item($r, 4).Value2,
        $ws.Cells.Item($r, 10).Value2,
        $ws.Cells.Item($r, 11).Value2,
        $ws.Cells.Item($r, 12).Value2,
        $ws.Cells.Item($r, 13).Value2,
        $ws.Cells.Item($r, 16).Value2
    )
}

foreach ($r in 2..46) {
    $src = $mapping[$r]
    $vals = $orig[$src]
    $ws.Cells.Item($r, 4).Value2 = $vals[0]
    $ws.Cells.Item($r, 10).Value2 = $vals[1]
    $ws.Cells.Item($r, 11).Value2 = $vals[2]
    $ws.Cells.Item($r, 12).Value2 = $vals[3]
    $ws.Cells.Item($r, 13).Value2 = $vals[4]
    $ws.Cells.Item($r, 16).Value2 = $vals[5]
}
